# Auto-generated edit script: update crypto price/volume figures
# D column values must stay TEXT (inlineStr) even when the new value looks
# numeric (e.g. "212.56") - Excel auto-converts numeric-looking strings to
# the Number type on plain .Value assignment, so we force text format first
# and restore the default "Normal" style afterwards (matches original: no
# explicit style index on these data cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.530.49"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.68%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.567.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.66%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.85%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "46.12"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.85%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "24.08"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.40%  "
$ws.Range("E10").Value = "  -1.74%  "
$ws.Range("E11").Value = "  -1.74%  "
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.791.83"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.565.39"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.81%  "
$ws.Range("E15").Value = "  -2.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.513.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "230.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.05%  "
$ws.Range("E20").Value = "  -1.88%  "
$ws.Range("E21").Value = "  -2.68%  "
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("E23").Value = "  -6.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.35%  "
$ws.Range("E25").Value = "  +8.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.28%  "
$ws.Range("E28").Value = "  -2.73%  "
$ws.Range("E29").Value = "  -3.59%  "
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0485"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.89%  "
$ws.Range("E32").Value = "  -3.60%  "
$ws.Range("E33").Value = "  -1.23%  "
$ws.Range("E34").Value = "  -1.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.392.34"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.44%  "
$ws.Range("E36").Value = "  +0.75%  "
$ws.Range("E37").Value = "  -3.92%  "
$ws.Range("E38").Value = "  +0.59%  "
$ws.Range("E39").Value = "  +2.98%  "
$ws.Range("E40").Value = "  -1.05%  "
$ws.Range("E41").Value = "  -3.71%  "
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("E43").Value = "  -0.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.787"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0463"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.48"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.971"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "62.74"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.704.67"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "86.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.48%  "
$ws.Range("E51").Value = "  -0.82%  "
